$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "About" (sheet1): reword the Notes paragraph, and drop the old
# internal "use China variables" reminder row (the blank spacer row
# above it collapses so the final note moves from row 13 up to row 12).
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Rows(12).Delete()

$wsAbout.Range("A10").Value = "The reserve margin (difference between the total generation available and the forecasted peak demand) in the U.S. "
$wsAbout.Range("A11").Value = "dataset doesn't vary by year, but the RM Reserve Margin variable is a time series to support countries that project "
$wsAbout.Range("A12").Value = "changes in future reserve margin by year."
$wsAbout.Range("A12").ClearFormats()

# ---------------------------------------------------------------------
# Sheet "RM" (sheet2): label the series as dimensionless, and correct
# the reserve-margin assumption from a flat 15% to 14.12%.
# ---------------------------------------------------------------------
$wsRM = $wb.Worksheets.Item("RM")

$wsRM.Range("A1:AK1").ClearFormats()
$wsRM.Range("A1").Value = "(dimensionless)"

$wsRM.Range("A2").ClearFormats()

$wsRM.Range("B2:AK2").Value = 0.1412
$wsRM.Range("B2:AK2").ClearFormats()
$wsRM.Range("B2:AK2").NumberFormat = "General"

# ---------------------------------------------------------------------
# Restore the on-disk cursor/selection state and make "About" the
# sheet that is active when the workbook is reopened.
# ---------------------------------------------------------------------
$wsRM.Range("A2").Select()

$wsAbout.Activate()
$wsAbout.Range("F17").Select()
